# v1.0.10 Add analyses for all three `MCP1`, `MCP1_pg_ml_2015`, and
# `MCP1_pg_ug_2015`. Add comparison between `MCP1`, `MCP1_pg_ml_2015`, and
# `MCP1_pg_ug_2015`. Add (and fixed) ordinal regression. Double checked
# which measurement to use.
#
# The "Results" sheet held one Cox-regression row per outcome, all run
# against the single CpG measure "MCP1_pg_ug_2015_rank". This edit expands
# the comparison to all three MCP1 measurements (MCP1_pg_ug_2015_rank,
# MCP1_pg_ml_2015_rank, MCP1_rank) crossed with all four outcomes
# (epmajor.3years, epstroke.3years, epcoronary.3years, epcvdeath.3years),
# refreshing the regression statistics to the re-checked values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Outcome x CpG-measure grid (rows 2-13), columns:
# A Dataset | B Outcome | C CpG | D Beta | E s.e.m. | F HR | G low95CI |
# H up95CI | I Z-value | J P-value | K SampleSize | L N_events
$rows = @(
    @("AEDB.CEA", "epmajor.3years", "MCP1_pg_ug_2015_rank", 0.256928881312197, 0.205011680561919, 1.29295317035662, 0.865113979442616, 1.9323787852929, 1.25324020859679, 0.210118300922242, 1029, 115),
    @("AEDB.CEA", "epmajor.3years", "MCP1_pg_ml_2015_rank", 0.201234418477267, 0.206883415752133, 1.22291141125608, 0.81525276200303, 1.83441552053863, 0.972694779548522, 0.330705020614303, 1029, 115),
    @("AEDB.CEA", "epmajor.3years", "MCP1_rank", -0.415016293511907, 0.278206701745752, 0.660329521530409, 0.382777488373779, 1.13913458928074, -1.49175519823093, 0.135763332809262, 493, 61),
    @("AEDB.CEA", "epstroke.3years", "MCP1_pg_ug_2015_rank", 0.36196934292784, 0.292494441463305, 1.43615491298066, 0.809515173025046, 2.54787186554027, 1.23752554447518, 0.215892036927472, 1029, 59),
    @("AEDB.CEA", "epstroke.3years", "MCP1_pg_ml_2015_rank", 0.237313571323117, 0.292023722247155, 1.2678386132819, 0.715300228600254, 2.24718892160017, 0.812651689722199, 0.416417786758067, 1029, 59),
    @("AEDB.CEA", "epstroke.3years", "MCP1_rank", -0.518645007996574, 0.399850914176228, 0.595326664576008, 0.271891334022138, 1.3035128126825, -1.29709596654314, 0.194598169172463, 493, 29),
    @("AEDB.CEA", "epcoronary.3years", "MCP1_pg_ug_2015_rank", -0.109643724433961, 0.244058288872441, 0.896153355972045, 0.555438175889808, 1.44586899547085, -0.449252205038885, 0.653249731750069, 1029, 78),
    @("AEDB.CEA", "epcoronary.3years", "MCP1_pg_ml_2015_rank", 0.250465326909837, 0.24976087104817, 1.28462304730306, 0.787362930181354, 2.09592845993155, 1.00282052131989, 0.31594746564528, 1029, 78),
    @("AEDB.CEA", "epcoronary.3years", "MCP1_rank", -0.00666446979910148, 0.331299977447844, 0.993357688527924, 0.518916432272481, 1.90157689367485, -0.0201161190847067, 0.983950741577137, 493, 42),
    @("AEDB.CEA", "epcvdeath.3years", "MCP1_pg_ug_2015_rank", -0.0313965789678052, 0.390911811274651, 0.969091175680942, 0.450415977363359, 2.08504527810089, -0.0803162709907128, 0.93598571764001, 1029, 33),
    @("AEDB.CEA", "epcvdeath.3years", "MCP1_pg_ml_2015_rank", 0.142597457350568, 0.400517430246023, 1.1532654696779, 0.526019651891625, 2.52846303131164, 0.356033087656073, 0.721815786217517, 1029, 33),
    @("AEDB.CEA", "epcvdeath.3years", "MCP1_rank", -0.262228624319768, 0.449374810239128, 0.769335114882511, 0.318859965440188, 1.85622713147509, -0.583540995945516, 0.559529167587422, 493, 23)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
